$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.365.05'
$ws.Range('E2').Value = '  +2.26%  '

$ws.Range('D3').Value = '1.825.57'
$ws.Range('E3').Value = '  +1.52%  '

$ws.Range('D4').Value = "'0.9999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').Value = "'314.16"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.63%  '

$ws.Range('E6').Value = '  -0.14%  '

$ws.Range('D7').Value = "'0.4694"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +5.23%  '

$ws.Range('D8').Value = "'0.3788"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.34%  '

$ws.Range('D9').Value = "'0.07442"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.63%  '

$ws.Range('D10').Value = "'0.8767"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.49%  '

$ws.Range('D11').Value = "'20.80"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.05%  '

$ws.Range('D12').Value = '1.826.35'
$ws.Range('E12').Value = '  -1.78%  '

$ws.Range('E13').Value = '  +1.33%  '

$ws.Range('D14').Value = "'5.432"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.79%  '

$ws.Range('D15').Value = "'93.08"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.08%  '

$ws.Range('E16').Value = '  +0.25%  '

$ws.Range('E17').Value = '  -0.10%  '

$ws.Range('D18').Value = "'0.000008805"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.31%  '

$ws.Range('D19').Value = "'0.9999"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.06%  '

$ws.Range('D20').Value = "'15.05"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.43%  '

$ws.Range('D21').Value = '27.362.97'
$ws.Range('E21').Value = '  +2.07%  '

$ws.Range('D22').Value = "'5.328"
$ws.Range('D22').Style = 'Normal'

$ws.Range('D23').Value = "'10.97"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.83%  '

$ws.Range('D24').Value = '2.051.94'
$ws.Range('E24').Value = '  -2.46%  '

$ws.Range('D25').Value = "'1.942"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.34%  '

$ws.Range('D26').Value = "'151.24"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.50%  '

$ws.Range('D27').Value = "'2.258"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.48%  '

$ws.Range('D28').Value = "'18.62"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.91%  '

$ws.Range('D29').Value = "'5.345"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.94%  '

$ws.Range('D30').Value = "'117.16"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.50%  '

$ws.Range('D31').Value = "'0.08999"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.24%  '

$ws.Range('E32').Value = '  +6.31%  '

$ws.Range('E33').Value = '  +3.23%  '

$ws.Range('D34').Value = "'4.547"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.19%  '

$ws.Range('E35').Value = '  -0.05%  '

$ws.Range('D36').Value = "'0.9995"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.16%  '

$ws.Range('E37').Value = '  +1.59%  '

$ws.Range('E38').Value = '  +0.84%  '

$ws.Range('D39').Value = "'0.05254"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.58%  '

$ws.Range('D40').Value = "'7.294"
$ws.Range('D40').Style = 'Normal'

$ws.Range('D41').Value = "'0.5354"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.17%  '

$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = "'2.371"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +20.80%  '

$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = "'2.895"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.69%  '

$ws.Range('E44').Value = '  +1.32%  '

$ws.Range('D45').Value = "'8.667"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.26%  '

$ws.Range('D46').Value = "'0.5091"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.18%  '

$ws.Range('D47').Value = "'10.62"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.79%  '

$ws.Range('D48').Value = "'105.83"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.31%  '

$ws.Range('D49').Value = "'1.684"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.45%  '

$ws.Range('D50').Value = "'0.9995"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.10%  '

$ws.Range('D51').Value = "'0.06383"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.21%  '

